$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# The document has a "different first page" header/footer, so there are
# two distinct header parts and two distinct footer parts (no separate
# even-page variant). Each one carries a single inline picture: the BTec
# logo in the headers, the Pearson logo in the footers. Rename every one
# of them, matching the commit's picture-name swap:
#   BTec logo:    image1.jpg -> image2.jpg
#   Pearson logo: image2.png -> image1.png
# (the pictures are identified by their alt text / description, since
# that's stable, unlike Name which Word reports blank until explicitly
# set in this session)

for ($h = 1; $h -le 3; $h++) {
    $hdr = $sec.Headers.Item($h)
    if ($hdr.Exists) {
        for ($i = 1; $i -le $hdr.Range.InlineShapes.Count; $i++) {
            $shp = $hdr.Range.InlineShapes.Item($i)
            if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                $shp.Name = "image2.jpg"
            }
        }
    }
}

for ($f = 1; $f -le 3; $f++) {
    $ftr = $sec.Footers.Item($f)
    if ($ftr.Exists) {
        for ($i = 1; $i -le $ftr.Range.InlineShapes.Count; $i++) {
            $shp = $ftr.Range.InlineShapes.Item($i)
            if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                $shp.Name = "image1.png"
            }
        }
    }
}
